# Auto-generated edit script applying the Halicarnassus_Profits.xlsx diff
# (workbook sheets correspond to the original "Sheets/Halicarnassus_Profits.xlsx" tabs)
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 157.4
$ws.Range("I5").Value = 163.77777
$ws.Range("K5").Value = 163.77777
$ws.Range("M5").Value = -48.77777

# Sheet ALC, row 62
$ws.Range("H62").Value = 8865.571
$ws.Range("J62").Value = 9911.799999999999
$ws.Range("L62").Value = 9911.799999999999
$ws.Range("N62").Value = -11159.8

# Sheet ALC, row 65
$ws.Range("H65").Value = 8865.571
$ws.Range("J65").Value = 9911.799999999999
$ws.Range("L65").Value = 49559
$ws.Range("N65").Value = -55799

# Sheet ALC, row 121
$ws.Range("H121").Value = 610.3
$ws.Range("J121").Value = 610.3
$ws.Range("L121").Value = 1830.9
$ws.Range("N121").Value = -5324.9

# Sheet ALC, row 132
$ws.Range("H132").Value = 19218.75
$ws.Range("I132").Value = 17884.691
$ws.Range("K132").Value = 53654.073
$ws.Range("M132").Value = -51124.073

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5770.722
$ws.Range("I61").Value = 5325.067
$ws.Range("K61").Value = 5325.067
$ws.Range("M61").Value = -5113.067

# Sheet ARM, row 74
$ws.Range("H74").Value = 4817.6875
$ws.Range("I74").Value = 3468.2307
$ws.Range("J74").Value = 10665.333
$ws.Range("K74").Value = 3468.2307
$ws.Range("L74").Value = 10665.333
$ws.Range("M74").Value = -2594.2307
$ws.Range("N74").Value = -12413.333

# Sheet ARM, row 77
$ws.Range("H77").Value = 4817.6875
$ws.Range("I77").Value = 3468.2307
$ws.Range("J77").Value = 10665.333
$ws.Range("K77").Value = 17341.1535
$ws.Range("L77").Value = 53326.665
$ws.Range("M77").Value = -12973.1535
$ws.Range("N77").Value = -62062.665

# Sheet ARM, row 82
$ws.Range("H82").Value = 28000
$ws.Range("J82").Value = 28000
$ws.Range("L82").Value = 28000
$ws.Range("N82").Value = -28722

# Sheet ARM, row 85
$ws.Range("H85").Value = 28000
$ws.Range("J85").Value = 28000
$ws.Range("L85").Value = 28000
$ws.Range("N85").Value = -30496

# Sheet ARM, row 88
$ws.Range("H88").Value = 2437.7273
$ws.Range("I88").Value = 2687.8572
$ws.Range("K88").Value = 2687.8572
$ws.Range("M88").Value = -2281.8572

# Sheet ARM, row 91
$ws.Range("H91").Value = 2437.7273
$ws.Range("I91").Value = 2687.8572
$ws.Range("K91").Value = 2687.8572
$ws.Range("M91").Value = -1283.8572

# Sheet ARM, row 122
$ws.Range("H122").Value = 1987.5
$ws.Range("I122").Value = 1816.6666
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5449.9998
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2999.9998
$ws.Range("N122").Value = -12400

# Sheet ARM, row 132
$ws.Range("H132").Value = 4939.8
$ws.Range("J132").Value = 4933
$ws.Range("L132").Value = 14799
$ws.Range("N132").Value = -19859

# Sheet ARM, row 136
$ws.Range("H136").Value = 5770.722
$ws.Range("I136").Value = 5325.067
$ws.Range("K136").Value = 15975.201
$ws.Range("M136").Value = -13425.201

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.625
$ws.Range("J7").Value = 264
$ws.Range("L7").Value = 264
$ws.Range("N7").Value = -490

# Sheet CRP, row 31
$ws.Range("H31").Value = 8903.308000000001
$ws.Range("I31").Value = 8356.857
$ws.Range("K31").Value = 8356.857
$ws.Range("M31").Value = -8061.857

# Sheet CRP, row 34
$ws.Range("H34").Value = 8903.308000000001
$ws.Range("I34").Value = 8356.857
$ws.Range("K34").Value = 8356.857
$ws.Range("M34").Value = -8154.857

# Sheet CRP, row 51
$ws.Range("H51").Value = 28187.5
$ws.Range("J51").Value = 42375
$ws.Range("L51").Value = 42375
$ws.Range("N51").Value = -43847

# Sheet CRP, row 59
$ws.Range("H59").Value = 53993.25
$ws.Range("J59").Value = 61001.4
$ws.Range("L59").Value = 61001.4
$ws.Range("N59").Value = -63291.4

# Sheet CRP, row 60
$ws.Range("H60").Value = 1000
$ws.Range("I60").Value = 1000
$ws.Range("K60").Value = 1000
$ws.Range("M60").Value = -489

# Sheet CRP, row 61
$ws.Range("H61").Value = 28187.5
$ws.Range("J61").Value = 42375
$ws.Range("L61").Value = 42375
$ws.Range("N61").Value = -43071

# Sheet CRP, row 122
$ws.Range("H122").Value = 1580.5
$ws.Range("J122").Value = 1896
$ws.Range("L122").Value = 5688
$ws.Range("N122").Value = -10588

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1067.8235
$ws.Range("I5").Value = 1067.8235
$ws.Range("K5").Value = 3203.4705
$ws.Range("M5").Value = -3091.4705

# Sheet CUL, row 60
$ws.Range("H60").Value = 791.087
$ws.Range("I60").Value = 273.1579
$ws.Range("J60").Value = 3251.25
$ws.Range("K60").Value = 819.4737
$ws.Range("L60").Value = 9753.75
$ws.Range("M60").Value = -568.4737
$ws.Range("N60").Value = -10255.75

# Sheet CUL, row 61
$ws.Range("H61").Value = 29
$ws.Range("I61").Value = 29
$ws.Range("K61").Value = 87
$ws.Range("M61").Value = 128

# Sheet CUL, row 121
$ws.Range("H121").Value = 562.5
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -5620

# Sheet CUL, row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Sheet CUL, row 135
$ws.Range("H135").Value = 1067.8235
$ws.Range("I135").Value = 1067.8235
$ws.Range("K135").Value = 9610.4115
$ws.Range("M135").Value = -7075.4115

# Sheet GSM, row 59
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 10000
$ws.Range("J59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("N59").Value = -11166

# Sheet GSM, row 80
$ws.Range("H80").Value = 4333
$ws.Range("J80").Value = 4499.5
$ws.Range("L80").Value = 4499.5
$ws.Range("N80").Value = -6495.5

# Sheet GSM, row 83
$ws.Range("H83").Value = 4333
$ws.Range("J83").Value = 4499.5
$ws.Range("L83").Value = 22497.5
$ws.Range("N83").Value = -32481.5

# Sheet GSM, row 132
$ws.Range("H132").Value = 71050.44500000001
$ws.Range("I132").Value = 77744.375
$ws.Range("K132").Value = 233233.125
$ws.Range("M132").Value = -230703.125

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4464.4
$ws.Range("I82").Value = 3927.75
$ws.Range("J82").Value = 4822.1665
$ws.Range("K82").Value = 3927.75
$ws.Range("L82").Value = 4822.1665
$ws.Range("M82").Value = -3566.75
$ws.Range("N82").Value = -5544.1665

# Sheet LTW, row 85
$ws.Range("H85").Value = 4464.4
$ws.Range("I85").Value = 3927.75
$ws.Range("J85").Value = 4822.1665
$ws.Range("K85").Value = 3927.75
$ws.Range("L85").Value = 4822.1665
$ws.Range("M85").Value = -2679.75
$ws.Range("N85").Value = -7318.1665

# Sheet LTW, row 132
$ws.Range("H132").Value = 7661.826
$ws.Range("I132").Value = 6790.0586
$ws.Range("K132").Value = 20370.1758
$ws.Range("M132").Value = -17840.1758

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Sheet WVR, row 132
$ws.Range("H132").Value = 2099.75
$ws.Range("I132").Value = 2099.75
$ws.Range("K132").Value = 6299.25
$ws.Range("M132").Value = -3769.25
